# edit load test function edit qa test
#
# QA/checkListUserInteract.xlsx - "чекЛист" sheet:
#   - row 9  (Ландшафтная ориентация устройства): результат not run -> failed,
#            комментарий added explaining the landscape-orientation bug,
#            row made taller so the comment is readable
#   - rows 10-12: результат not run -> passed
#   - row 13: new test case "Плавность скроллинга" / failed, with comment
#   - conditional formatting (pass/fail/blocked/skipped colouring) extended
#     to include the new row
#   - page set up as portrait / A4
#   - selection left on D14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("чекЛист")

# --- Row 9: результат -> failed, add комментарий, grow the row to fit it
$ws.Range("C9").Value = "failed"
$ws.Range("D9").Value = "Приложение переключилось на ландшафтную ориентацию. При ландшафтной ориентации приложение не выполняет основные функции. Ландшафтная ориентация не предусмотрена ТЗ."
$ws.Range("D9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 68

# --- Rows 10-12: результат -> passed
$ws.Range("C10").Value = "passed"
$ws.Range("C11").Value = "passed"
$ws.Range("C12").Value = "passed"

# --- Row 13: new test case "Плавность скроллинга" (failed)
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Плавность скроллинга"
$ws.Range("C13").ClearFormats()
$ws.Range("C13").Value = "failed"
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("D13").Value = "При скроллинге более 10 элементов видно затормаживание"

# --- Conditional formatting on the result column now covers row 13 too
$fcs = $ws.Range("C9:C12").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("C9:C13"))
}

# --- Page setup: portrait / A4
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ends on D14
$ws.Range("D14").Select() | Out-Null
